$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 2) beneath the header row, mirroring the style of
# the added row in the target workbook:
# A=identifier, B=alternativeIdentifiers, C=title, D=date_s,
# E=levelOfDescription, F=extentAndMedium, G=notes, H=file_path
$ws.Range("A2").Value = "MCH205-1"
$ws.Range("C2").Value = "LOCAL AUTHOROTIES AGAINST APARTHEID"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24C | GRAP COUNT NUMER: NONE"

$ws.Range("A2:H2").Font.Size = 10
$ws.Range("A2:H2").Font.Name = "Calibri"
$ws.Range("A2:H2").Font.Color = $ws.Range("A1").Font.Color

$ws.Range("A2").Select()
